$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rsquo = [char]0x2019

# Row 4: question 3 - "how old are you?"
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "how old are you?"
$ws.Range("C4").Value = "not important"
$ws.Range("D4").Value = "it${rsquo}s none of your buisiness"
$ws.Range("E4").Value = "who the hell would have cared?"
$ws.Range("F4").Value = "get lost"
$ws.Range("G4").Value = "it${rsquo}s not important"
$ws.Range("H4").Value = "none"

# Row 5: question 4 - "how big is earth"
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "how big is earth"
$ws.Range("C5").Value = "none"
$ws.Range("D5").Value = "none"
$ws.Range("E5").Value = "none"
$ws.Range("F5").Value = "none"
$ws.Range("G5").Value = 85000
$ws.Range("H5").Value = "https://www.esri.com/news/arcuser/0610/graphics/nospin_1-lg.jpg"

# Row 6: question 5 - "who is an old guy"
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "who is an old guy"
$ws.Range("C6").Value = "https://upload.wikimedia.org/wikipedia/commons/5/53/Random_Old_Guy.jpg"
$ws.Range("D6").Value = "https://jooinn.com/images/happy-young-man-1.png"
$ws.Range("E6").Value = "https://www.stockvault.net/data/2018/04/25/250732/preview16.jpg"
$ws.Range("F6").Value = "https://d2v9y0dukr6mq2.cloudfront.net/video/thumbnail/NGyZeGzFlijx95hou/4k-close-up-face-of-a-young-man-without-emotions-beautiful-frowning-guy-in-a-white-shirt-looking-to-the-camera_bjv1icufl_thumbnail-full01.png"
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = "none"

[void]$ws.Range("C6").Select()
